$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.852.39'
$ws.Range('E2').Value = '  -1.19%  '
$ws.Range('D3').Value = '3.033.43'
$ws.Range('E3').Value = '  -0.59%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').Value = '''538.47'
$ws.Range('E5').Value = '  +0.76%  '
$ws.Range('D6').Value = '''136.17'
$ws.Range('E6').Value = '  +2.61%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = '3.026.13'
$ws.Range('E8').Value = '  -0.68%  '
$ws.Range('D9').Value = '''0.500'
$ws.Range('E9').Value = '  +1.33%  '
$ws.Range('D10').Value = '''0.151'
$ws.Range('E10').Value = '  -1.16%  '
$ws.Range('D11').Value = '''6.18'
$ws.Range('E11').Value = '  +0.94%  '
$ws.Range('D12').Value = '''0.453'
$ws.Range('E12').Value = '  +0.53%  '
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').Value = '''35.11'
$ws.Range('E13').Value = '  +3.51%  '
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').Value = '''0.0000224'
$ws.Range('E14').Value = '  +1.03%  '
$ws.Range('D15').Value = '3.513.28'
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('D17').Value = '61.718.90'
$ws.Range('E17').Value = '  -1.49%  '
$ws.Range('D18').Value = '3.026.33'
$ws.Range('E18').Value = '  -1.06%  '
$ws.Range('D19').Value = '''6.70'
$ws.Range('E19').Value = '  +1.71%  '
$ws.Range('D20').Value = '''472.10'
$ws.Range('E20').Value = '  -1.75%  '
$ws.Range('D21').Value = '''13.43'
$ws.Range('E21').Value = '  +2.00%  '
$ws.Range('D22').Value = '''0.687'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').Value = '''7.08'
$ws.Range('E23').Value = '  +0.23%  '
$ws.Range('D24').Value = '''79.93'
$ws.Range('E24').Value = '  +1.38%  '
$ws.Range('D25').Value = '''12.30'
$ws.Range('E25').Value = '  +2.51%  '
$ws.Range('D26').Value = '''0.999'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').Value = '''2.71'
$ws.Range('E27').Value = '  +1.09%  '
$ws.Range('D28').Value = '''7.96'
$ws.Range('E28').Value = '  -0.69%  '
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('D30').Value = '''1.93'
$ws.Range('E30').Value = '  +4.82%  '
$ws.Range('D31').Value = '''25.90'
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('E32').Value = '  +4.07%  '
$ws.Range('D33').Value = '''5.58'
$ws.Range('E33').Value = '  +5.60%  '
$ws.Range('B34').Value = 'Stacks'
$ws.Range('C34').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D34').Value = '''2.32'
$ws.Range('E34').Value = '  -0.15%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = '''55.94'
$ws.Range('E35').Value = '  -1.02%  '
$ws.Range('D36').Value = '''5.98'
$ws.Range('E36').Value = '  +1.11%  '
$ws.Range('D37').Value = '''475.96'
$ws.Range('E37').Value = '  +0.85%  '
$ws.Range('D38').Value = '3.224.09'
$ws.Range('E38').Value = '  +4.77%  '
$ws.Range('D39').Value = '''0.0800'
$ws.Range('E39').Value = '  +1.40%  '
$ws.Range('D40').Value = '''0.0393'
$ws.Range('E40').Value = '  +0.17%  '
$ws.Range('D41').Value = '''0.119'
$ws.Range('E41').Value = '  +4.44%  '
$ws.Range('D42').Value = '''8.24'
$ws.Range('E42').Value = '  +2.67%  '
$ws.Range('D43').Value = '''2.55'
$ws.Range('E43').Value = '  -2.75%  '
$ws.Range('D44').Value = '''27.48'
$ws.Range('E44').Value = '  +13.16%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').Value = '''0.253'
$ws.Range('E45').Value = '  +1.27%  '
$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').Value = '''1.00'
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').Value = '''2.04'
$ws.Range('E47').Value = '  +2.30%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').Value = '''0.110'
$ws.Range('E48').Value = '  +2.28%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').Value = '''119.79'
$ws.Range('E49').Value = '  -0.92%  '
$ws.Range('D50').Value = '0.0₃0507'
$ws.Range('E50').Value = '  -5.22%  '
$ws.Range('E51').Value = '  +8.34%  '
